$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "52.086.67"
Set-TextValue 2 5 "  +0.88%  "

# Row 3
Set-TextValue 3 4 "2.899.10"
Set-TextValue 3 5 "  +3.48%  "

# Row 4
Set-TextValue 4 5 "  -0.08%  "

# Row 5
Set-TextValue 5 4 "350.83"
Set-TextValue 5 5 "  -1.21%  "

# Row 6
Set-TextValue 6 4 "111.62"
Set-TextValue 6 5 "  +1.97%  "

# Row 7
Set-TextValue 7 5 "  +0.33%  "

# Row 8
Set-TextValue 8 5 "  +0.00%  "

# Row 9
Set-TextValue 9 4 "0.622"
Set-TextValue 9 5 "  -0.02%  "

# Row 10
Set-TextValue 10 4 "39.88"
Set-TextValue 10 5 "  -0.75%  "

# Row 11
Set-TextValue 11 5 "  +0.39%  "

# Row 12
Set-TextValue 12 4 "0.0858"
Set-TextValue 12 5 "  +2.38%  "

# Row 13
Set-TextValue 13 4 "19.96"
Set-TextValue 13 5 "  -0.65%  "

# Row 14
Set-TextValue 14 4 "7.78"
Set-TextValue 14 5 "  -0.30%  "

# Row 15
Set-TextValue 15 4 "3.351.78"
Set-TextValue 15 5 "  +3.41%  "

# Row 16
Set-TextValue 16 4 "0.999"
Set-TextValue 16 5 "  +5.95%  "

# Row 17
Set-TextValue 17 4 "2.898.89"
Set-TextValue 17 5 "  +3.52%  "

# Row 18
Set-TextValue 18 4 "52.061.86"
Set-TextValue 18 5 "  +0.89%  "

# Row 19
Set-TextValue 19 2 "InternetComputer(DFINITY)"
Set-TextValue 19 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 19 4 "14.62"
Set-TextValue 19 5 "  +9.16%  "

# Row 20
Set-TextValue 20 2 "Uniswap"
Set-TextValue 20 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue 20 4 "7.70"
Set-TextValue 20 5 "  -0.95%  "

# Row 21
Set-TextValue 21 2 "ImmutableX"
Set-TextValue 21 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 21 4 "3.33"
Set-TextValue 21 5 "  +4.35%  "

# Row 22
Set-TextValue 22 4 "0.0₃0979"
Set-TextValue 22 5 "  +0.76%  "

# Row 23
Set-TextValue 23 4 "70.76"
Set-TextValue 23 5 "  +0.48%  "

# Row 24
Set-TextValue 24 4 "269.69"
Set-TextValue 24 5 "  +0.54%  "

# Row 25
Set-TextValue 25 5 "  +0.42%  "

# Row 26
Set-TextValue 26 4 "26.60"
Set-TextValue 26 5 "  +1.95%  "

# Row 27
Set-TextValue 27 5 "  +0.11%  "

# Row 28
Set-TextValue 28 4 "0.165"
Set-TextValue 28 5 "  -0.02%  "

# Row 29
Set-TextValue 29 4 "10.53"
Set-TextValue 29 5 "  +1.58%  "

# Row 30
Set-TextValue 30 4 "38.34"
Set-TextValue 30 5 "  +2.46%  "

# Row 31
Set-TextValue 31 4 "2.24"
Set-TextValue 31 5 "  +0.27%  "

# Row 32
Set-TextValue 32 4 "6.46"
Set-TextValue 32 5 "  +1.80%  "

# Row 33
Set-TextValue 33 5 "  +7.76%  "

# Row 34
Set-TextValue 34 4 "0.0943"
Set-TextValue 34 5 "  +10.34%  "

# Row 35
Set-TextValue 35 4 "52.90"
Set-TextValue 35 5 "  +1.36%  "

# Row 36
Set-TextValue 36 4 "0.0458"
Set-TextValue 36 5 "  +2.39%  "

# Row 37
Set-TextValue 37 5 "  -0.22%  "

# Row 38
Set-TextValue 38 4 "3.29"
Set-TextValue 38 5 "  +4.52%  "

# Row 39
Set-TextValue 39 4 "18.63"
Set-TextValue 39 5 "  -1.63%  "

# Row 40
Set-TextValue 40 4 "2.04"
Set-TextValue 40 5 "  +2.76%  "

# Row 41
Set-TextValue 41 4 "2.65"
Set-TextValue 41 5 "  +6.48%  "

# Row 42
Set-TextValue 42 5 "  +1.90%  "

# Row 43
Set-TextValue 43 4 "22.69"
Set-TextValue 43 5 "  +3.39%  "

# Row 44
Set-TextValue 44 4 "121.88"
Set-TextValue 44 5 "  +1.66%  "

# Row 45
Set-TextValue 45 5 "  +0.57%  "

# Row 46
Set-TextValue 46 4 "3.56"
Set-TextValue 46 5 "  +4.02%  "

# Row 47
Set-TextValue 47 2 "Maker"
Set-TextValue 47 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 47 4 "2.197.36"
Set-TextValue 47 5 "  +2.75%  "

# Row 48
Set-TextValue 48 2 "ApeXProtocol"
Set-TextValue 48 3 "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue 48 4 "2.53"
Set-TextValue 48 5 "  +6.48%  "

# Row 49
Set-TextValue 49 4 "0.272"
Set-TextValue 49 5 "  +23.09%  "

# Row 50
Set-TextValue 50 4 "0.946"
Set-TextValue 50 5 "  +2.79%  "

# Row 51
Set-TextValue 51 4 "0.0324"
Set-TextValue 51 5 "  +10.56%  "
